# Applies the "EC" (Estado de Cuenta) data update described by the commit:
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# (Delete previous EC rows and add new ones; the underlying database changes.)
#
# The worker table occupies rows 16-20, columns C (document number),
# D (worker name) and E (period). Column B (doc type "CC") and the
# amount columns F/G stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New contents for the table body (rows 16-20): DocNumber, Name, Period
$rows = @(
    @{ Row = 16; Doc = "73196852"; Name = "JORGE LUIS TAPIAS ROJAS"; Period = "1712" },
    @{ Row = 17; Doc = "73196852"; Name = "JORGE LUIS TAPIAS ROJAS"; Period = "1711" },
    @{ Row = 18; Doc = "9153002";  Name = "FELIX VALENCIA PEREZ";    Period = "1801" },
    @{ Row = 19; Doc = "9153002";  Name = "FELIX VALENCIA PEREZ";    Period = "1712" },
    @{ Row = 20; Doc = "9153002";  Name = "FELIX VALENCIA PEREZ";    Period = "1711" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc      # Column C - N Doc Trabajador
    $ws.Cells.Item($r.Row, 4).Value = $r.Name     # Column D - Nombre Trabajador
    $ws.Cells.Item($r.Row, 5).Value = $r.Period   # Column E - Periodo Mora
}

# The new text (e.g. "73196852" / "JORGE LUIS TAPIAS ROJAS") is wider than
# the values it replaces, so Excel's "best fit" columns grow to keep
# fitting the data without wrapping/truncating.
$ws.Columns.Item(2).ColumnWidth = 17.666666666666664
$ws.Columns.Item(3).ColumnWidth = 15.833333333333332
$ws.Columns.Item(5).ColumnWidth = 12.666666666666668
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333336

$wb.Save()
